$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 6: SingleUseId3 / Default / Left / LTR / 00000
$ws.Range("B6").Value = "SingleUseId3"
$ws.Range("C6").Value = "Default"
$ws.Range("D6").Value = "Left"
$ws.Range("E6").Value = "LTR"

# Row 7: SingleUseId4 / Default / Left / LTR / 00000
$ws.Range("B7").Value = "SingleUseId4"
$ws.Range("C7").Value = "Default"
$ws.Range("D7").Value = "Left"
$ws.Range("E7").Value = "LTR"

# F6/F7 need the literal text "00000" (leading zeros preserved) without
# altering the cell's number format (stay General/style 0), so build it
# via a formula and convert to a static value with Paste Special (values).
$ws.Range("F6").Formula = '=TEXT(0,"00000")'
$ws.Range("F7").Formula = '=TEXT(0,"00000")'
$ws.Range("F6:F7").Copy()
$ws.Range("F6:F7").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$wb.Save()
